# The source record (row 5) and the record below it (row 6) traded places:
# row 5's data becomes row 6's data and vice versa. Only the columns that
# actually differ between the two rows need to be swapped; the rest of the
# row (Lokalnamn, Noggrannhet, Lan, Kommun, etc.) is identical already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

foreach ($col in $cols) {
    $cell5 = $ws.Range($col + "5")
    $cell6 = $ws.Range($col + "6")

    $val5 = $cell5.Value()
    $val6 = $cell6.Value()

    $cell5.Value = $val6
    $cell6.Value = $val5
}
